# Update column G ("K") values on Sheet1 (rows 2-30) per regenerated save_data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @(
    3,
    7,
    7,
    7,
    3,
    4,
    7,
    3,
    6,
    2,
    7,
    2,
    4,
    3,
    8,
    2,
    10,
    6,
    4,
    10,
    3,
    3,
    6,
    5,
    5,
    3,
    3,
    3,
    2
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newValues[$i]
}
